# Move column-header definitions out of the sheet and into short field
# names (the long descriptions now live in the README instead).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) cell-by-cell. The order below matches the
# order new shared-string entries were written by the original commit.
$ws.Range("A1").Value = "DOI"
$ws.Range("C1").Value = "Digestion"
$ws.Range("D1").Value = "Filtration"
$ws.Range("H1").Value = "Controls"
$ws.Range("E1").Value = "Filter_Size"
$ws.Range("F1").Value = "Microplastic_Identification_Method"
$ws.Range("G1").Value = "Spectral_Analysis"
$ws.Range("B1").Value = "Sample_device_and_deployment_methods"

# Update the view state: scrolled/selected cell moved from A54 to B6.
$ws.Range("B6").Select()
